$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.8057775497436523
$ws.Range("E2").Value = 4579.600213179509
$ws.Range("K2").Value = 0.1075019016274153
$ws.Range("L2").Value = 0.1074888211835606
$ws.Range("M2").Value = 0.1065843255760504
$ws.Range("N2").Value = 0.1065843255760504
$ws.Range("O2").Value = 0.1060407155884529
$ws.Range("P2").Value = 0.1060407155884529
$ws.Range("Q2").Value = 0.1060407155884529
$ws.Range("R2").Value = 0.1060407155884529
$ws.Range("S2").Value = 0.1053953846563404
$ws.Range("T2").Value = 0.1053953846563404
$ws.Range("U2").Value = 0.1053953846563404
$ws.Range("V2").Value = 0.1053066515441534
$ws.Range("W2").Value = 0.1052889596742081
$ws.Range("X2").Value = 0.1052889596742081
$ws.Range("Y2").Value = 0.105270959321238
$ws.Range("C3").Value = 0.9119932651519775
$ws.Range("E3").Value = 4579.307614346204
$ws.Range("L3").Value = 0.1071457251935886
$ws.Range("M3").Value = 0.1071457251935886
$ws.Range("N3").Value = 0.1060023192380948
$ws.Range("O3").Value = 0.1060023192380948
$ws.Range("P3").Value = 0.1060023192380948
$ws.Range("Q3").Value = 0.1060023192380948
$ws.Range("R3").Value = 0.1060023192380948
$ws.Range("S3").Value = 0.1060023192380948
$ws.Range("T3").Value = 0.1053555684563418
$ws.Range("U3").Value = 0.1052899878317202
$ws.Range("V3").Value = 0.1052899878317202
$ws.Range("W3").Value = 0.1052899878317202
$ws.Range("X3").Value = 0.1052652556402769
$ws.Range("Y3").Value = 0.1052652556402769
$ws.Range("C4").Value = 0.7844038009643555
$ws.Range("E4").Value = 4578.891474699201
$ws.Range("K4").Value = 0.1070619727165266
$ws.Range("L4").Value = 0.1070619727165266
$ws.Range("M4").Value = 0.1070438236302132
$ws.Range("N4").Value = 0.1056391856218288
$ws.Range("O4").Value = 0.1056391856218288
$ws.Range("P4").Value = 0.1056391856218288
$ws.Range("Q4").Value = 0.105347232971681
$ws.Range("R4").Value = 0.105347232971681
$ws.Range("S4").Value = 0.105347232971681
$ws.Range("T4").Value = 0.105347232971681
$ws.Range("U4").Value = 0.105347232971681
$ws.Range("V4").Value = 0.1053204914032758
$ws.Range("W4").Value = 0.105272245909678
$ws.Range("X4").Value = 0.1052571437563197
$ws.Range("Y4").Value = 0.1052571437563197
$ws.Range("C5").Value = 0.7968897819519043
$ws.Range("E5").Value = 4582.737127651533
$ws.Range("J5").Value = 0.1065624492353679
$ws.Range("K5").Value = 0.1055772594239685
$ws.Range("L5").Value = 0.1055772594239685
$ws.Range("M5").Value = 0.1055772594239685
$ws.Range("N5").Value = 0.1055772594239685
$ws.Range("O5").Value = 0.1055772594239685
$ws.Range("P5").Value = 0.1055772594239685
$ws.Range("Q5").Value = 0.1055772594239685
$ws.Range("R5").Value = 0.1055202967545155
$ws.Range("S5").Value = 0.1053586914978238
$ws.Range("T5").Value = 0.1053586914978238
$ws.Range("U5").Value = 0.1053586914978238
$ws.Range("V5").Value = 0.1053321077514919
$ws.Range("W5").Value = 0.1053321077514919
$ws.Range("X5").Value = 0.1053321077514919
$ws.Range("Y5").Value = 0.1053321077514919
$ws.Range("C6").Value = 0.7968709468841553
$ws.Range("E6").Value = 4581.088744419015
$ws.Range("K6").Value = 0.1069630472723423
$ws.Range("L6").Value = 0.1069630472723423
$ws.Range("M6").Value = 0.1066792259690654
$ws.Range("N6").Value = 0.1059793579356022
$ws.Range("O6").Value = 0.1054106989892522
$ws.Range("P6").Value = 0.1054106989892522
$ws.Range("Q6").Value = 0.1054106989892522
$ws.Range("R6").Value = 0.1054106989892522
$ws.Range("S6").Value = 0.1054106989892522
$ws.Range("T6").Value = 0.1053125804047961
$ws.Range("U6").Value = 0.1053125804047961
$ws.Range("V6").Value = 0.1053125804047961
$ws.Range("W6").Value = 0.1053125804047961
$ws.Range("X6").Value = 0.1052999755247371
$ws.Range("Y6").Value = 0.1052999755247371
$ws.Range("C7").Value = 0.7968652248382568
$ws.Range("E7").Value = 4578.457432586334
$ws.Range("J7").Value = 0.1063000828664116
$ws.Range("K7").Value = 0.1063000828664116
$ws.Range("L7").Value = 0.1053915303629414
$ws.Range("M7").Value = 0.1053915303629414
$ws.Range("N7").Value = 0.1053915303629414
$ws.Range("O7").Value = 0.1053915303629414
$ws.Range("P7").Value = 0.1053915303629414
$ws.Range("Q7").Value = 0.1052750472760308
$ws.Range("R7").Value = 0.1052750472760308
$ws.Range("S7").Value = 0.1052750472760308
$ws.Range("T7").Value = 0.1052750472760308
$ws.Range("U7").Value = 0.1052750472760308
$ws.Range("V7").Value = 0.1052750472760308
$ws.Range("W7").Value = 0.1052750472760308
$ws.Range("X7").Value = 0.1052486828964198
$ws.Range("Y7").Value = 0.1052486828964198
$ws.Range("C8").Value = 0.8124985694885254
$ws.Range("E8").Value = 4578.439842065173
$ws.Range("K8").Value = 0.1075019016274153
$ws.Range("L8").Value = 0.1055091129974938
$ws.Range("M8").Value = 0.1053020261018551
$ws.Range("N8").Value = 0.1053020261018551
$ws.Range("O8").Value = 0.1053020261018551
$ws.Range("P8").Value = 0.1053020261018551
$ws.Range("Q8").Value = 0.1053020261018551
$ws.Range("R8").Value = 0.1053020261018551
$ws.Range("S8").Value = 0.1053020261018551
$ws.Range("T8").Value = 0.1053020261018551
$ws.Range("U8").Value = 0.1053020261018551
$ws.Range("V8").Value = 0.1052561916668784
$ws.Range("W8").Value = 0.1052561916668784
$ws.Range("X8").Value = 0.1052483400012704
$ws.Range("Y8").Value = 0.1052483400012704
$ws.Range("C9").Value = 0.7812516689300537
$ws.Range("E9").Value = 4578.291388007952
$ws.Range("J9").Value = 0.1065073651198012
$ws.Range("K9").Value = 0.1065073651198012
$ws.Range("L9").Value = 0.1060608090164503
$ws.Range("M9").Value = 0.1060608090164503
$ws.Range("N9").Value = 0.1053788354323011
$ws.Range("O9").Value = 0.1053788354323011
$ws.Range("P9").Value = 0.1053788354323011
$ws.Range("Q9").Value = 0.1053788354323011
$ws.Range("R9").Value = 0.1053788354323011
$ws.Range("S9").Value = 0.1053788354323011
$ws.Range("T9").Value = 0.1053788354323011
$ws.Range("U9").Value = 0.1052640622691247
$ws.Range("V9").Value = 0.1052640622691247
$ws.Range("W9").Value = 0.1052634987310408
$ws.Range("X9").Value = 0.1052634987310408
$ws.Range("Y9").Value = 0.105245446159999
$ws.Range("C10").Value = 0.781224250793457
$ws.Range("E10").Value = 4579.043172884622
$ws.Range("J10").Value = 0.1068047185099058
$ws.Range("K10").Value = 0.1068047185099058
$ws.Range("L10").Value = 0.1063103485318554
$ws.Range("M10").Value = 0.1063103485318554
$ws.Range("N10").Value = 0.1060852370270839
$ws.Range("O10").Value = 0.1060852370270839
$ws.Range("P10").Value = 0.1060852370270839
$ws.Range("Q10").Value = 0.1053582414450147
$ws.Range("R10").Value = 0.1053582414450147
$ws.Range("S10").Value = 0.1053582414450147
$ws.Range("T10").Value = 0.1053582414450147
$ws.Range("U10").Value = 0.1053582414450147
$ws.Range("V10").Value = 0.1053582414450147
$ws.Range("W10").Value = 0.1052601008359575
$ws.Range("X10").Value = 0.1052601008359575
$ws.Range("Y10").Value = 0.1052601008359575
$ws.Range("C11").Value = 0.7968995571136475
$ws.Range("E11").Value = 4581.900914603236
$ws.Range("J11").Value = 0.1075019016274153
$ws.Range("K11").Value = 0.1075019016274153
$ws.Range("L11").Value = 0.1071568231277742
$ws.Range("M11").Value = 0.1071568231277742
$ws.Range("N11").Value = 0.1065465200169519
$ws.Range("O11").Value = 0.1055378804817546
$ws.Range("P11").Value = 0.1055378804817546
$ws.Range("Q11").Value = 0.1055378804817546
$ws.Range("R11").Value = 0.1055378804817546
$ws.Range("S11").Value = 0.1055378804817546
$ws.Range("T11").Value = 0.1053718667400415
$ws.Range("U11").Value = 0.1052627496367762
$ws.Range("V11").Value = 0.1052627496367762
$ws.Range("W11").Value = 0.1053158073022073
$ws.Range("X11").Value = 0.1053158073022073
$ws.Range("Y11").Value = 0.1052818065065667
